$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two team labels (B4/B5) to match reordering in shared strings table
$ws.Range("B4").Value = "Kurni_s2l"
$ws.Range("B5").Value = "EDN_S2l"

# Update the numeric data (C2:D8); formulas in E recalc automatically
$ws.Range("C2").Value = 2462401585
$ws.Range("D2").Value = 256164284

$ws.Range("C3").Value = 1658736781
$ws.Range("D3").Value = 171937370

$ws.Range("C4").Value = 5966177823
$ws.Range("D4").Value = 601636347

$ws.Range("C5").Value = 1629709461
$ws.Range("D5").Value = 163243379

$ws.Range("C6").Value = 6063602699
$ws.Range("D6").Value = 534260929

$ws.Range("C7").Value = 5934351328
$ws.Range("D7").Value = 506386477

$ws.Range("C8").Value = 3202314817
$ws.Range("D8").Value = 273083942

# Update the active cell selection to C1
$ws.Range("C1").Select()

$wb.Save()
